$d = $word.ActiveDocument

$d.Content.Find.Execute("371÷5=74, 1", $true, $false, $false, $false, $false, $true, 1, $false, "258÷2=129, 0", 2) | Out-Null
$d.Content.Find.Execute("846÷4=211, 2", $true, $false, $false, $false, $false, $true, 1, $false, "475÷9=52, 7", 2) | Out-Null
$d.Content.Find.Execute("135÷2=67, 1", $true, $false, $false, $false, $false, $true, 1, $false, "427÷7=61, 0", 2) | Out-Null
$d.Content.Find.Execute("690÷7=98, 4", $true, $false, $false, $false, $false, $true, 1, $false, "129÷3=43, 0", 2) | Out-Null
$d.Content.Find.Execute("996÷7=142, 2", $true, $false, $false, $false, $false, $true, 1, $false, "854÷9=94, 8", 2) | Out-Null
$d.Content.Find.Execute("389÷5=77, 4", $true, $false, $false, $false, $false, $true, 1, $false, "329÷3=109, 2", 2) | Out-Null
$d.Content.Find.Execute("469÷6=78, 1", $true, $false, $false, $false, $false, $true, 1, $false, "212÷4=53, 0", 2) | Out-Null
$d.Content.Find.Execute("862÷2=431, 0", $true, $false, $false, $false, $false, $true, 1, $false, "753÷4=188, 1", 2) | Out-Null
$d.Content.Find.Execute("225÷6=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "222÷6=37, 0", 2) | Out-Null
$d.Content.Find.Execute("322÷4=80, 2", $true, $false, $false, $false, $false, $true, 1, $false, "127÷7=18, 1", 2) | Out-Null
$d.Content.Find.Execute("639÷9=71, 0", $true, $false, $false, $false, $false, $true, 1, $false, "545÷8=68, 1", 2) | Out-Null
$d.Content.Find.Execute("732÷5=146, 2", $true, $false, $false, $false, $false, $true, 1, $false, "963÷3=321, 0", 2) | Out-Null
$d.Content.Find.Execute("931÷5=186, 1", $true, $false, $false, $false, $false, $true, 1, $false, "736÷2=368, 0", 2) | Out-Null
$d.Content.Find.Execute("107÷5=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "666÷5=133, 1", 2) | Out-Null
$d.Content.Find.Execute("205÷3=68, 1", $true, $false, $false, $false, $false, $true, 1, $false, "654÷5=130, 4", 2) | Out-Null
$d.Content.Find.Execute("876÷3=292, 0", $true, $false, $false, $false, $false, $true, 1, $false, "983÷9=109, 2", 2) | Out-Null
$d.Content.Find.Execute("419÷4=104, 3", $true, $false, $false, $false, $false, $true, 1, $false, "898÷9=99, 7", 2) | Out-Null
$d.Content.Find.Execute("103÷3=34, 1", $true, $false, $false, $false, $false, $true, 1, $false, "449÷9=49, 8", 2) | Out-Null
$d.Content.Find.Execute("148÷3=49, 1", $true, $false, $false, $false, $false, $true, 1, $false, "470÷8=58, 6", 2) | Out-Null
$d.Content.Find.Execute("635÷8=79, 3", $true, $false, $false, $false, $false, $true, 1, $false, "133÷9=14, 7", 2) | Out-Null
$d.Content.Find.Execute("273÷8=34, 1", $true, $false, $false, $false, $false, $true, 1, $false, "988÷6=164, 4", 2) | Out-Null
$d.Content.Find.Execute("804÷9=89, 3", $true, $false, $false, $false, $false, $true, 1, $false, "315÷5=63, 0", 2) | Out-Null
$d.Content.Find.Execute("174÷5=34, 4", $true, $false, $false, $false, $false, $true, 1, $false, "298÷5=59, 3", 2) | Out-Null
$d.Content.Find.Execute("658÷5=131, 3", $true, $false, $false, $false, $false, $true, 1, $false, "183÷7=26, 1", 2) | Out-Null
$d.Content.Find.Execute("269÷4=67, 1", $true, $false, $false, $false, $false, $true, 1, $false, "317÷2=158, 1", 2) | Out-Null
